$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24 (shifts existing rows 24-49 down to 25-50,
# carries the date-format style of column D along with it)
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new weekly price entry
$ws.Range("A24").Value = 11
$ws.Range("B24").Value = "Vega Monumental Concepción"
$ws.Range("C24").Value = "Bíobío"
$ws.Range("D24").Value = "2023-04-28"
$ws.Range("E24").Value = 8
$ws.Range("F24").Value = 100114007
$ws.Range("G24").Value = "Jengibre"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 50
$ws.Range("K24").Value = 20000
$ws.Range("L24").Value = 20000
$ws.Range("M24").Value = 20000
$ws.Range("N24").Value = '$/caja 13 kilos'
$ws.Range("O24").Value = "Perú"
$ws.Range("P24").Value = 1538
$ws.Range("Q24").Value = 13
$ws.Range("R24").Value = "Hortaliza"
